# Update computed market-price / profit figures on several sheets
# (scheduled runner refresh of Sheets). Plain value overwrites — the
# source cells hold no formulas, so Range.Value assignment is sufficient.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H2").Value = 5041
$ws.Range("J2").Value = 5331
$ws.Range("L2").Value = 5331
$ws.Range("N2").Value = -5557

$ws.Range("H33").Value = 1167.8889
$ws.Range("I33").Value = 1257.625
$ws.Range("K33").Value = 1257.625
$ws.Range("M33").Value = -1028.625

$ws.Range("H98").Value = 5364
$ws.Range("J98").Value = 6333.222
$ws.Range("L98").Value = 6333.222
$ws.Range("N98").Value = -9329.222

$ws.Range("H122").Value = 5364
$ws.Range("J122").Value = 6333.222
$ws.Range("L122").Value = 18999.666
$ws.Range("N122").Value = -23899.666

$ws.Range("H131").Value = 8471.227999999999
$ws.Range("I131").Value = 1965.6364
$ws.Range("K131").Value = 5896.9092
$ws.Range("M131").Value = -856.9092000000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 4843.25
$ws.Range("I32").Value = 4299.6333
$ws.Range("J32").Value = 12997.5
$ws.Range("K32").Value = 4299.6333
$ws.Range("L32").Value = 12997.5
$ws.Range("M32").Value = -4012.6333
$ws.Range("N32").Value = -13571.5

$ws.Range("H45").Value = 1680.2
$ws.Range("I45").Value = 1658.7858
$ws.Range("K45").Value = 1658.7858
$ws.Range("M45").Value = -1281.7858

$ws.Range("H97").Value = 448.33334
$ws.Range("I97").Value = 467
$ws.Range("J97").Value = 388.6
$ws.Range("K97").Value = 467
$ws.Range("L97").Value = 388.6
$ws.Range("M97").Value = 29
$ws.Range("N97").Value = -1380.6

$ws.Range("H122").Value = 2839.8
$ws.Range("J122").Value = 3466.6667
$ws.Range("L122").Value = 10400.0001
$ws.Range("N122").Value = -15300.0001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H22").Value = 1065.8
$ws.Range("I22").Value = 1137.1428
$ws.Range("K22").Value = 1137.1428
$ws.Range("M22").Value = -964.1428000000001

$ws.Range("H86").Value = 2690.8823
$ws.Range("I86").Value = 1207.5
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 1207.5
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -84.5
$ws.Range("N86").Value = -5746

$ws.Range("H89").Value = 2690.8823
$ws.Range("I89").Value = 1207.5
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 6037.5
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -421.5
$ws.Range("N89").Value = -28732

$ws.Range("H134").Value = 850
$ws.Range("I134").Value = 850
$ws.Range("K134").Value = 2550
$ws.Range("M134").Value = -15

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H22").Value = 14119.389
$ws.Range("I22").Value = 237.14285
$ws.Range("J22").Value = 62707.25
$ws.Range("K22").Value = 237.14285
$ws.Range("L22").Value = 62707.25
$ws.Range("M22").Value = 112.85715
$ws.Range("N22").Value = -63407.25

$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 1000
$ws.Range("M58").Value = -797

$ws.Range("H62").Value = 3347.4
$ws.Range("J62").Value = 3347.4
$ws.Range("L62").Value = 3347.4
$ws.Range("N62").Value = -4595.4

$ws.Range("H65").Value = 3347.4
$ws.Range("J65").Value = 3347.4
$ws.Range("L65").Value = 16737
$ws.Range("N65").Value = -22977

$ws.Range("H99").Value = 10040.611
$ws.Range("J99").Value = 7659.4287
$ws.Range("L99").Value = 7659.4287
$ws.Range("N99").Value = -10655.4287

$ws.Range("H126").Value = 10040.611
$ws.Range("J126").Value = 7659.4287
$ws.Range("L126").Value = 22978.2861
$ws.Range("N126").Value = -27918.2861

$ws.Range("H132").Value = 1631.15
$ws.Range("I132").Value = 1352.4
$ws.Range("K132").Value = 4057.2
$ws.Range("M132").Value = -1527.2

$ws.Range("H134").Value = 3546.2
$ws.Range("I134").Value = 3353.8
$ws.Range("J134").Value = 3738.6
$ws.Range("K134").Value = 10061.4
$ws.Range("L134").Value = 11215.8
$ws.Range("M134").Value = -7526.400000000001
$ws.Range("N134").Value = -16285.8

$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

$ws.Range("H141").Value = 320833.16
$ws.Range("J141").Value = 320833.16
$ws.Range("L141").Value = 320833.16
$ws.Range("N141").Value = -331193.16

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 18334446
$ws.Range("I4").Value = 18334446
$ws.Range("K4").Value = 55003338
$ws.Range("M4").Value = -55003226

$ws.Range("H109").Value = 1849.5714
$ws.Range("I109").Value = 1849.5714
$ws.Range("K109").Value = 5548.7142
$ws.Range("M109").Value = -4508.7142

$ws.Range("H122").Value = 792.6
$ws.Range("J122").Value = 791
$ws.Range("L122").Value = 7119
$ws.Range("N122").Value = -12019

$ws.Range("H131").Value = 2215.3142
$ws.Range("J131").Value = 2818.0952
$ws.Range("L131").Value = 8454.285600000001
$ws.Range("N131").Value = -18534.2856

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H5").Value = 144000
$ws.Range("I5").Value = 144000
$ws.Range("K5").Value = 144000
$ws.Range("M5").Value = -143888

$ws.Range("H70").Value = 5683.778
$ws.Range("I70").Value = 5394
$ws.Range("K70").Value = 5394
$ws.Range("M70").Value = -5124

$ws.Range("H73").Value = 5683.778
$ws.Range("I73").Value = 5394
$ws.Range("K73").Value = 5394
$ws.Range("M73").Value = -4458

$ws.Range("H95").Value = 45729.668
$ws.Range("J95").Value = 45729.668
$ws.Range("L95").Value = 45729.668
$ws.Range("N95").Value = -51221.668

$ws.Range("H97").Value = 689.8182
$ws.Range("I97").Value = 620.5714
$ws.Range("K97").Value = 620.5714
$ws.Range("M97").Value = -124.5714

# Row 126: leve reverted to a 0-profit placeholder; the HQ-profit cell
# (N126) is removed entirely rather than merely zeroed.
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H16").Value = 1098
$ws.Range("I16").Value = 994.6
$ws.Range("J16").Value = 1227.25
$ws.Range("K16").Value = 994.6
$ws.Range("L16").Value = 1227.25
$ws.Range("M16").Value = -824.6
$ws.Range("N16").Value = -1567.25

$ws.Range("H122").Value = 8472.546
$ws.Range("I122").Value = 8760
$ws.Range("J122").Value = 7856.5713
$ws.Range("K122").Value = 26280
$ws.Range("L122").Value = 23569.7139
$ws.Range("M122").Value = -23830
$ws.Range("N122").Value = -28469.7139

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 122 previously had no NQ-profit (M122) cell; it is introduced here.
$ws.Range("H122").Value = 4900
$ws.Range("I122").Value = 2200
$ws.Range("K122").Value = 6600
$ws.Range("M122").Value = -4150
